# Add certificate file names + hyperlinks to column E (CERTIFICATE) for rows 2-11,
# widen column E to fit the new content, and move the active selection to H11 -
# matching the "Updated Excel with Test Data" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$certificates = @(
    "TEST1@USER1.pfx",
    "TEST2@USER2.pfx",
    "TEST3@USER3.pfx",
    "TEST4@USER4.pfx",
    "TEST5@USER5.pfx",
    "TEST6@USER6.pfx",
    "TEST7@USER7.pfx",
    "TEST8@USER8.pfx",
    "TEST9@USER9.pfx",
    "TEST10@USER10.pfx"
)

for ($i = 0; $i -lt $certificates.Count; $i++) {
    $row = $i + 2
    $cell = $ws.Range("E" + $row)
    $cell.Value = $certificates[$i]
    $ws.Hyperlinks.Add($cell, "mailto:" + $certificates[$i])
}

# Widen column E (CERTIFICATE) so the new, longer values fit.
$ws.Columns("E").ColumnWidth = 17.5

# Leave the selection where the author ended up after entering the data.
[void]$ws.Range("H11").Select()
